# Regenerate save_data to use K instead of Strike# (column G "K" values),
# recompute std/mean, calc and write s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values per row (row number => new value), per the diff.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    21 = 0
    22 = 3
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
